# New weekly price observation for Cilantro / Terminal La Palmera de La Serena.
# A new row is inserted at row 226 (pushing the existing rows 226-247 down to
# 227-248) and the newly inserted row is populated with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift rows 226:247 down by inserting a new row at 226.
$ws.Rows.Item(226).Insert()

# Populate the newly inserted row 226 with the new observation.
$ws.Cells.Item(226, 1).Value = 8
$ws.Cells.Item(226, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(226, 3).Value = "Coquimbo"
$ws.Cells.Item(226, 4).Value = 45132
$ws.Cells.Item(226, 5).Value = 4
$ws.Cells.Item(226, 6).Value = 100112040
$ws.Cells.Item(226, 7).Value = "Cilantro"
$ws.Cells.Item(226, 8).Value = "Sin especificar"
$ws.Cells.Item(226, 9).Value = "Primera"
$ws.Cells.Item(226, 10).Value = 2000
$ws.Cells.Item(226, 11).Value = 2500
$ws.Cells.Item(226, 12).Value = 3000
$ws.Cells.Item(226, 13).Value = 2750
$ws.Cells.Item(226, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(226, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(226, 16).Value = 1833
$ws.Cells.Item(226, 17).Value = 1.5
$ws.Cells.Item(226, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D.
$ws.Cells.Item(226, 4).NumberFormat = $ws.Cells.Item(227, 4).NumberFormat
